$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N; this shifts old N:P (Late, Date/heading, Outstanding) to O:Q
$ws.Columns("N").Insert()

# Give the new blank column the same width as the column to its left ("In Advance")
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# --- Make "Repayment schedule" the active sheet/tab and update its selection ---
$ws.Activate() | Out-Null
$ws.Range("L15").Select() | Out-Null
